$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell H1 with same style as G1 (bold header style)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H69 with the "Save" indicator values
$saveValues = @(1,0,0,1,0,1,0,1,0,1,0,0,0,0,0,0,0,1,0,0,0,0,1,0,0,0,0,1,0,1,1,1,0,0,0,0,1,1,1,1,0,0,0,0,1,0,0,0,1,0,1,0,1,1,0,1,0,0,1,0,0,1,0,0,0,0,1,1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

